$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 0.0003064172267913819
$ws.Range("M2").Value = 0.001490354537963867
$ws.Range("N2").Value = 0.0001442432403564453

$ws.Range("L3").Value = 0.00009489059448242188
$ws.Range("M3").Value = 0.0004427433013916016
$ws.Range("N3").Value = 0.00006103515625

$ws.Range("F4").Value = 0.65
$ws.Range("G4").Value = 0.6666666666666666
$ws.Range("H4").Value = 81
$ws.Range("I4").Value = 14
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 26
$ws.Range("L4").Value = 0.0003180978298187256
$ws.Range("M4").Value = 0.001538753509521484
$ws.Range("N4").Value = 0.0001435279846191406

$ws.Range("F5").Value = 0.6578947368421053
$ws.Range("G5").Value = 0.6410256410256411
$ws.Range("H5").Value = 82
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 14
$ws.Range("K5").Value = 25
$ws.Range("L5").Value = 0.00009478831291198731
$ws.Range("M5").Value = 0.0006144046783447266
$ws.Range("N5").Value = 0.00006771087646484375

$ws.Range("E6").Value = 0.8805970149253731
$ws.Range("F6").Value = 0.7948717948717948
$ws.Range("H6").Value = 87
$ws.Range("I6").Value = 8
$ws.Range("L6").Value = 0.01129823732376099
$ws.Range("M6").Value = 0.0503842830657959
$ws.Range("N6").Value = 0.008116245269775391

$ws.Range("E7").Value = 0.8805970149253731
$ws.Range("F7").Value = 0.7948717948717948
$ws.Range("H7").Value = 87
$ws.Range("I7").Value = 8
$ws.Range("L7").Value = 0.007691164493560791
$ws.Range("M7").Value = 0.01374673843383789
$ws.Range("N7").Value = 0.006586313247680664

$ws.Range("L8").Value = 0.0003345158100128174
$ws.Range("M8").Value = 0.001645803451538086
$ws.Range("N8").Value = 0.0001373291015625

$ws.Range("L9").Value = 0.00009734678268432618
$ws.Range("M9").Value = 0.0007107257843017578
$ws.Range("N9").Value = 0.00007677078247070312

$ws.Range("L10").Value = 0.0003338503837585449
$ws.Range("M10").Value = 0.001607418060302734
$ws.Range("N10").Value = 0.00014495849609375

$ws.Range("L11").Value = 0.0001149258613586426
$ws.Range("M11").Value = 0.0004529953002929688
$ws.Range("N11").Value = 0.00008845329284667969

$ws.Range("L12").Value = 0.001303379774093628
$ws.Range("M12").Value = 0.004743576049804688
$ws.Range("N12").Value = 0.0005819797515869141

$ws.Range("L13").Value = 0.0005609757900238037
$ws.Range("M13").Value = 0.001481294631958008
$ws.Range("N13").Value = 0.0004479885101318359

$ws.Range("L14").Value = 0.0002442803382873535
$ws.Range("M14").Value = 0.0009565353393554688
$ws.Range("N14").Value = 0.0001006126403808594

$ws.Range("L15").Value = 0.00007474923133850098
$ws.Range("M15").Value = 0.0004110336303710938

$ws.Range("L16").Value = 0.0004940464496612549
$ws.Range("M16").Value = 0.004242897033691406
$ws.Range("N16").Value = 0.0001595020294189453

$ws.Range("L17").Value = 0.000166006326675415
$ws.Range("M17").Value = 0.007096529006958008
$ws.Range("N17").Value = 0.0001032352447509766
